$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

$ws.Range("A2").Value = "0d94751"
$ws.Range("M2").Value = "japanese"
$ws.Range("N2").Value = "BIRDS"
$ws.Range("A3").Value = "19f2049"
$ws.Range("M3").Value = "english"
$ws.Range("N3").Value = "DOGS"
$ws.Range("A4").Formula = "=""5311690"""
$ws.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("N4").Value = "CATS"
$ws.Range("A5").Value = "a9b5545"
$ws.Range("N5").Value = "CATS"
$ws.Range("A6").Value = "58dd74d"
$ws.Range("M6").Value = "japanese"
$ws.Range("N6").Value = "CATS"
$ws.Range("A7").Value = "07bdb8d"
$ws.Range("M7").Value = "japanese"
$ws.Range("N7").Value = "BIRDS"
$ws.Range("A8").Value = "cfed5dc"
$ws.Range("M8").Value = "japanese"
$ws.Range("N8").Value = "CATS"
$ws.Range("A9").Value = "24aea3a"
$ws.Range("N9").Value = "REPTILES"
$ws.Range("A10").Value = "0b026af"
$ws.Range("M10").Value = "english"
$ws.Range("A11").Value = "ca0e083"
$ws.Range("N11").Value = "DOGS"
$ws.Range("A12").Value = "1697c6a"
$ws.Range("M12").Value = "english"
$ws.Range("N12").Value = "DOGS"
$ws.Range("A13").Value = "c73e5a2"
$ws.Range("M13").Value = "japanese"
$ws.Range("N13").Value = "FISH"
$ws.Range("A14").Value = "1fa0325"
$ws.Range("M14").Value = "japanese"
$ws.Range("N14").Value = "CATS"
$ws.Range("A15").Value = "0614d2b"
$ws.Range("M15").Value = "japanese"
$ws.Range("N15").Value = "BIRDS"
$ws.Range("A16").Value = "88a29a6"
$ws.Range("N16").Value = "BIRDS"
$ws.Range("A17").Value = "e1af405"
$ws.Range("M17").Value = "japanese"
$ws.Range("N17").Value = "REPTILES"
$ws.Range("A18").Value = "012ed68"
$ws.Range("N18").Value = "REPTILES"
$ws.Range("A19").Value = "8fa1a0f"
$ws.Range("M19").Value = "english"
$ws.Range("N19").Value = "CATS"
$ws.Range("A20").Value = "231f1d5"
$ws.Range("N20").Value = "REPTILES"
$ws.Range("A21").Value = "b296f2b"
$ws.Range("N21").Value = "CATS"
$ws.Range("A22").Value = "6870f83"
$ws.Range("N22").Value = "DOGS"
$ws.Range("A23").Value = "d2ef047"
$ws.Range("A24").Value = "8a8e186"
$ws.Range("M24").Value = "japanese"
$ws.Range("N24").Value = "BIRDS"
$ws.Range("A25").Value = "e91ea49"
$ws.Range("N25").Value = "DOGS"
$ws.Range("A26").Value = "96a135d"
$ws.Range("M26").Value = "english"
$ws.Range("N26").Value = "DOGS"
$ws.Range("A27").Value = "6f0f461"
$ws.Range("N27").Value = "BIRDS"
$ws.Range("A28").Value = "a51c04c"
$ws.Range("N28").Value = "REPTILES"
$ws.Range("A29").Value = "c25733f"
$ws.Range("M29").Value = "japanese"
$ws.Range("N29").Value = "BIRDS"
$ws.Range("A30").Value = "e2ddb80"
$ws.Range("N30").Value = "FISH"
$ws.Range("A31").Value = "60f3311"
$ws.Range("M31").Value = "japanese"
$ws.Range("N31").Value = "REPTILES"
$ws.Range("A32").Value = "f658364"
$ws.Range("N32").Value = "FISH"
$ws.Range("A33").Value = "4c845b3"
$ws.Range("N33").Value = "REPTILES"
$ws.Range("A34").Value = "3556f67"
$ws.Range("A35").Value = "306e33c"
$ws.Range("A36").Value = "d3e37a4"
$ws.Range("A37").Value = "e73b883"
$ws.Range("N37").Value = "CATS"
$ws.Range("A38").Value = "883875d"
$ws.Range("M38").Value = "english"
$ws.Range("N38").Value = "CATS"
$ws.Range("A39").Value = "e640589"
$ws.Range("M39").Value = "japanese"
$ws.Range("N39").Value = "FISH"
$ws.Range("A40").Value = "e543986"
$ws.Range("A41").Value = "89aefd2"
$ws.Range("M41").Value = "english"
$ws.Range("N41").Value = "BIRDS"
$ws.Range("A42").Value = "5fb78ae"
$ws.Range("M42").Value = "english"
$ws.Range("N42").Value = "BIRDS"
$ws.Range("A43").Value = "1a65664"
$ws.Range("N43").Value = "REPTILES"
$ws.Range("A44").Value = "138adf6"
$ws.Range("M44").Value = "english"
$ws.Range("N44").Value = "BIRDS"
$ws.Range("A45").Value = "43d0f8a"
$ws.Range("M45").Value = "japanese"
$ws.Range("A46").Value = "586c386"
$ws.Range("M46").Value = "english"
$ws.Range("N46").Value = "CATS"
$ws.Range("A47").Formula = "=""79e0726"""
$ws.Range("A47").Copy()
$ws.Range("A47").PasteSpecial(-4163)
$ws.Range("N47").Value = "DOGS"
$ws.Range("A48").Value = "ca8645a"
$ws.Range("M48").Value = "japanese"
$ws.Range("N48").Value = "FISH"
$ws.Range("A49").Value = "2d656ae"
$ws.Range("M49").Value = "english"
$ws.Range("N49").Value = "CATS"
$ws.Range("A50").Value = "89b83b5"
$ws.Range("N50").Value = "CATS"
$ws.Range("A51").Value = "d687cdc"
$ws.Range("M51").Value = "japanese"

$excel.CutCopyMode = $false
